$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.193778157234192
$ws.Range("B1").Value = 2.32668137550354
$ws.Range("C1").Value = 3.498147487640381
$ws.Range("D1").Value = 3.344908952713013
$ws.Range("E1").Value = 1.141029477119446
